$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("C2").Value = 23602299
$ws.Range("D2").Value = 95.90000000000001
$ws.Range("E2").Value = 22642378
$ws.Range("F2").Value = 0.3
$ws.Range("G2").Value = 75499
$ws.Range("H2").Value = 1.2
$ws.Range("I2").Value = 286415
$ws.Range("J2").Value = 2.5
$ws.Range("K2").Value = 598007
$ws.Range("N2").Value = 23602299

# Update row 3 values (ocap -> non_pdi)
$ws.Range("A3").Value = "non_pdi (5-17 y.o.)"
$ws.Range("B3").Value = "non_pdi"
$ws.Range("C3").Value = 12525253
$ws.Range("D3").Value = 94.40000000000001
$ws.Range("E3").Value = 11822828
$ws.Range("F3").Value = 0.4
$ws.Range("G3").Value = 52304
$ws.Range("H3").Value = 1.8
$ws.Range("I3").Value = 227095
$ws.Range("J3").Value = 3.4
$ws.Range("K3").Value = 423026
$ws.Range("N3").Value = 12525253

# Update row 4 values (idp -> pdi)
$ws.Range("A4").Value = "pdi (5-17 y.o.)"
$ws.Range("B4").Value = "pdi"
$ws.Range("C4").Value = 11077047
$ws.Range("D4").Value = 97.7
$ws.Range("E4").Value = 10819550
$ws.Range("F4").Value = 0.2
$ws.Range("G4").Value = 23195
$ws.Range("H4").Value = 0.5
$ws.Range("I4").Value = 59320
$ws.Range("J4").Value = 1.6
$ws.Range("K4").Value = 174981
$ws.Range("N4").Value = 11077047

# Delete rows 5 and 6 entirely (ret and ndsp rows removed)
$ws.Range("A5:N6").EntireRow.Delete()
